$wb = $excel.ActiveWorkbook

# Map of cell -> new value to apply on both the "展览" and "全部类型" sheets
$updates = @{
    "F2" = 2213
    "F3" = 632
    "F4" = 1599
    "F5" = 7424
    "F7" = 191
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
